# vm_pu.xlsx - Case_2_121 (res_bus): refresh bus voltage magnitude results
# for the run with the 380 kV slack-bus setpoint (commit: "case with 380 kV done").
# Slack bus voltage (column B) drops from 1.05 to 1.02 p.u. and every other bus
# voltage (columns C:F, I:N) is updated to the recomputed load-flow solution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033615620209729
$ws.Cells.Item(2, 4).Value = 1.0379528171458
$ws.Cells.Item(2, 5).Value = 1.053801346495265
$ws.Cells.Item(2, 6).Value = 1.0596881405365
$ws.Cells.Item(2, 9).Value = 1.038247871030924
$ws.Cells.Item(2, 10).Value = 1.038739372639384
$ws.Cells.Item(2, 11).Value = 1.040742233260011
$ws.Cells.Item(2, 12).Value = 1.056546312603732
$ws.Cells.Item(2, 13).Value = 1.062416962913989
$ws.Cells.Item(2, 14).Value = 1.040214501693603

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034438901845615
$ws.Cells.Item(3, 4).Value = 1.038573976036072
$ws.Cells.Item(3, 5).Value = 1.054850870582612
$ws.Cells.Item(3, 6).Value = 1.060752876648986
$ws.Cells.Item(3, 9).Value = 1.038441229757854
$ws.Cells.Item(3, 10).Value = 1.039206081832141
$ws.Cells.Item(3, 11).Value = 1.041173805061382
$ws.Cells.Item(3, 12).Value = 1.05740834848967
$ws.Cells.Item(3, 13).Value = 1.063295344189091
$ws.Cells.Item(3, 14).Value = 1.040681873666946

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034972012839618
$ws.Cells.Item(4, 4).Value = 1.038976211702422
$ws.Cells.Item(4, 5).Value = 1.05553121278048
$ws.Cells.Item(4, 6).Value = 1.061442932080747
$ws.Cells.Item(4, 9).Value = 1.03856537903234
$ws.Cells.Item(4, 10).Value = 1.03950780001215
$ws.Cells.Item(4, 11).Value = 1.041452673092877
$ws.Cells.Item(4, 12).Value = 1.057966768169873
$ws.Cells.Item(4, 13).Value = 1.063864212402863
$ws.Cells.Item(4, 14).Value = 1.040984020321379

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035196225322524
$ws.Cells.Item(5, 4).Value = 1.039145382881274
$ws.Cells.Item(5, 5).Value = 1.055817522160584
$ws.Cells.Item(5, 6).Value = 1.061733293715592
$ws.Cells.Item(5, 9).Value = 1.03861733935213
$ws.Cells.Item(5, 10).Value = 1.039634575733708
$ws.Cells.Item(5, 11).Value = 1.041569815306685
$ws.Cells.Item(5, 12).Value = 1.058201676550986
$ws.Cells.Item(5, 13).Value = 1.064103482601734
$ws.Cells.Item(5, 14).Value = 1.041110976079002

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035233876975742
$ws.Cells.Item(6, 4).Value = 1.039173791602768
$ws.Cells.Item(6, 5).Value = 1.055865611950644
$ws.Cells.Item(6, 6).Value = 1.061782062065943
$ws.Cells.Item(6, 9).Value = 1.038626050094602
$ws.Cells.Item(6, 10).Value = 1.03965585799959
$ws.Cells.Item(6, 11).Value = 1.041589478461979
$ws.Cells.Item(6, 12).Value = 1.058241127400293
$ws.Cells.Item(6, 13).Value = 1.064143664013699
$ws.Cells.Item(6, 14).Value = 1.041132288568143

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034975008415595
$ws.Cells.Item(7, 4).Value = 1.038978471897766
$ws.Cells.Item(7, 5).Value = 1.055535037308245
$ws.Cells.Item(7, 6).Value = 1.061446810877584
$ws.Cells.Item(7, 9).Value = 1.038566074241152
$ws.Cells.Item(7, 10).Value = 1.039509494258538
$ws.Cells.Item(7, 11).Value = 1.041454238723201
$ws.Cells.Item(7, 12).Value = 1.057969906443717
$ws.Cells.Item(7, 13).Value = 1.063867409080459
$ws.Cells.Item(7, 14).Value = 1.040985716973791

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033893769935827
$ws.Cells.Item(8, 4).Value = 1.038162676755167
$ws.Cells.Item(8, 5).Value = 1.054155782848347
$ws.Cells.Item(8, 6).Value = 1.06004774472924
$ws.Cells.Item(8, 9).Value = 1.03831341733325
$ws.Cells.Item(8, 10).Value = 1.0388971552816
$ws.Cells.Item(8, 11).Value = 1.040888164620031
$ws.Cells.Item(8, 12).Value = 1.056837511790429
$ws.Cells.Item(8, 13).Value = 1.062713712552403
$ws.Cells.Item(8, 14).Value = 1.040372508405269

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031991559196709
$ws.Cells.Item(9, 4).Value = 1.036727540806749
$ws.Cells.Item(9, 5).Value = 1.051734821215178
$ws.Cells.Item(9, 6).Value = 1.057590880334606
$ws.Cells.Item(9, 9).Value = 1.037860825667902
$ws.Cells.Item(9, 10).Value = 1.037816083932692
$ws.Cells.Item(9, 11).Value = 1.03988774707793
$ws.Cells.Item(9, 12).Value = 1.054846910103996
$ws.Cells.Item(9, 13).Value = 1.060684597963257
$ws.Cells.Item(9, 14).Value = 1.039289901811045

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.030725565283513
$ws.Cells.Item(10, 4).Value = 1.035772488914496
$ws.Cells.Item(10, 5).Value = 1.050127262074469
$ws.Cells.Item(10, 6).Value = 1.055958723519588
$ws.Cells.Item(10, 9).Value = 1.037554172618341
$ws.Cells.Item(10, 10).Value = 1.037094051494449
$ws.Cells.Item(10, 11).Value = 1.039218900556164
$ws.Cells.Item(10, 12).Value = 1.053523133498391
$ws.Cells.Item(10, 13).Value = 1.059334495440541
$ws.Cells.Item(10, 14).Value = 1.038566844003921

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.030177902340665
$ws.Cells.Item(11, 4).Value = 1.035359364808684
$ws.Cells.Item(11, 5).Value = 1.049432703897379
$ws.Cells.Item(11, 6).Value = 1.055253358085909
$ws.Cells.Item(11, 9).Value = 1.037420227609329
$ws.Cells.Item(11, 10).Value = 1.036781103642953
$ws.Cells.Item(11, 11).Value = 1.038928844931198
$ws.Cells.Item(11, 12).Value = 1.052950714214966
$ws.Cells.Item(11, 13).Value = 1.058750523634074
$ws.Cells.Item(11, 14).Value = 1.038253451730581

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.02997455554525
$ws.Cells.Item(12, 4).Value = 1.035205976732175
$ws.Cells.Item(12, 5).Value = 1.049174944107635
$ws.Cells.Item(12, 6).Value = 1.054991560591461
$ws.Cells.Item(12, 9).Value = 1.037370300429971
$ws.Cells.Item(12, 10).Value = 1.036664816363061
$ws.Cells.Item(12, 11).Value = 1.03882104020782
$ws.Cells.Item(12, 12).Value = 1.05273821086292
$ws.Cells.Item(12, 13).Value = 1.058533706236577
$ws.Cells.Item(12, 14).Value = 1.038136999309413

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.030018170498841
$ws.Cells.Item(13, 4).Value = 1.035238876046471
$ws.Cells.Item(13, 5).Value = 1.049230224033938
$ws.Cells.Item(13, 6).Value = 1.055047707677218
$ws.Cells.Item(13, 9).Value = 1.037381017845882
$ws.Cells.Item(13, 10).Value = 1.036689762377393
$ws.Cells.Item(13, 11).Value = 1.038844167618139
$ws.Cells.Item(13, 12).Value = 1.052783788171538
$ws.Cells.Item(13, 13).Value = 1.058580209958395
$ws.Cells.Item(13, 14).Value = 1.038161980749947

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.030161091999271
$ws.Cells.Item(14, 4).Value = 1.035346684375375
$ws.Cells.Item(14, 5).Value = 1.049411392695027
$ws.Cells.Item(14, 6).Value = 1.055231713602722
$ws.Cells.Item(14, 9).Value = 1.037416104160146
$ws.Cells.Item(14, 10).Value = 1.036771492202321
$ws.Cells.Item(14, 11).Value = 1.038919935082059
$ws.Cells.Item(14, 12).Value = 1.052933146196607
$ws.Cells.Item(14, 13).Value = 1.058732599482456
$ws.Cells.Item(14, 14).Value = 1.0382438266406

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.030249161242187
$ws.Cells.Item(15, 4).Value = 1.035413117253126
$ws.Cells.Item(15, 5).Value = 1.049523047194775
$ws.Cells.Item(15, 6).Value = 1.055345113135094
$ws.Cells.Item(15, 9).Value = 1.03743769895159
$ws.Cells.Item(15, 10).Value = 1.036821842776728
$ws.Cells.Item(15, 11).Value = 1.038966609315858
$ws.Cells.Item(15, 12).Value = 1.053025186351347
$ws.Cells.Item(15, 13).Value = 1.058826504409015
$ws.Cells.Item(15, 14).Value = 1.038294248718598

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.030761922957995
$ws.Cells.Item(16, 4).Value = 1.03579991557596
$ws.Cells.Item(16, 5).Value = 1.050173389826737
$ws.Cells.Item(16, 6).Value = 1.056005565227114
$ws.Cells.Item(16, 9).Value = 1.037563037669574
$ws.Cells.Item(16, 10).Value = 1.037114814514551
$ws.Cells.Item(16, 11).Value = 1.039238141398478
$ws.Cells.Item(16, 12).Value = 1.053561139679034
$ws.Cells.Item(16, 13).Value = 1.059373265104016
$ws.Cells.Item(16, 14).Value = 1.038587636509894

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.031083705224996
$ws.Cells.Item(17, 4).Value = 1.036042657403797
$ws.Cells.Item(17, 5).Value = 1.050581741604139
$ws.Cells.Item(17, 6).Value = 1.056420216607677
$ws.Cells.Item(17, 9).Value = 1.037641348580173
$ws.Cells.Item(17, 10).Value = 1.037298507528878
$ws.Cells.Item(17, 11).Value = 1.039408349066235
$ws.Cells.Item(17, 12).Value = 1.053897539669594
$ws.Cells.Item(17, 13).Value = 1.059716403315483
$ws.Cells.Item(17, 14).Value = 1.038771590389371

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.031271445532817
$ws.Cells.Item(18, 4).Value = 1.036184284992257
$ws.Cells.Item(18, 5).Value = 1.050820073422315
$ws.Cells.Item(18, 6).Value = 1.056662207792008
$ws.Cells.Item(18, 9).Value = 1.037686913834152
$ws.Cells.Item(18, 10).Value = 1.037405623249153
$ws.Cells.Item(18, 11).Value = 1.039507585730342
$ws.Cells.Item(18, 12).Value = 1.054093831755583
$ws.Cells.Item(18, 13).Value = 1.059916610901033
$ws.Cells.Item(18, 14).Value = 1.038878858226255

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.031335468591893
$ws.Cells.Item(19, 4).Value = 1.036232583142781
$ws.Cells.Item(19, 5).Value = 1.050901363344271
$ws.Cells.Item(19, 6).Value = 1.056744742902194
$ws.Cells.Item(19, 9).Value = 1.037702431364047
$ws.Cells.Item(19, 10).Value = 1.037442141939789
$ws.Cells.Item(19, 11).Value = 1.039541415612378
$ws.Cells.Item(19, 12).Value = 1.054160775094951
$ws.Cells.Item(19, 13).Value = 1.059984886812871
$ws.Cells.Item(19, 14).Value = 1.03891542877762

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.031049175835363
$ws.Cells.Item(20, 4).Value = 1.036016609314251
$ws.Cells.Item(20, 5).Value = 1.05053791409799
$ws.Cells.Item(20, 6).Value = 1.056375714773586
$ws.Cells.Item(20, 9).Value = 1.037632958159286
$ws.Cells.Item(20, 10).Value = 1.037278802030936
$ws.Cells.Item(20, 11).Value = 1.039390091791846
$ws.Cells.Item(20, 12).Value = 1.053861439302566
$ws.Cells.Item(20, 13).Value = 1.059679581530295
$ws.Cells.Item(20, 14).Value = 1.038751856907361

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.030119002971116
$ws.Cells.Item(21, 4).Value = 1.03531493571347
$ws.Cells.Item(21, 5).Value = 1.049358036682708
$ws.Cells.Item(21, 6).Value = 1.05517752273234
$ws.Cells.Item(21, 9).Value = 1.037405776914838
$ws.Cells.Item(21, 10).Value = 1.036747426026589
$ws.Cells.Item(21, 11).Value = 1.038897625242356
$ws.Cells.Item(21, 12).Value = 1.052889160699222
$ws.Cells.Item(21, 13).Value = 1.058687721916215
$ws.Cells.Item(21, 14).Value = 1.038219726288139

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029534627628143
$ws.Cells.Item(22, 4).Value = 1.034874140259909
$ws.Cells.Item(22, 5).Value = 1.048617532301847
$ws.Cells.Item(22, 6).Value = 1.054425367822161
$ws.Cells.Item(22, 9).Value = 1.037261932396927
$ws.Cells.Item(22, 10).Value = 1.036413071801786
$ws.Cells.Item(22, 11).Value = 1.038587615564302
$ws.Cells.Item(22, 12).Value = 1.052278537444797
$ws.Cells.Item(22, 13).Value = 1.058064654719771
$ws.Cells.Item(22, 14).Value = 1.037884897241985

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.029844371875824
$ws.Cells.Item(23, 4).Value = 1.035107778242807
$ws.Cells.Item(23, 5).Value = 1.049009961200727
$ws.Cells.Item(23, 6).Value = 1.054823985696363
$ws.Cells.Item(23, 9).Value = 1.037338282325511
$ws.Cells.Item(23, 10).Value = 1.036590343328051
$ws.Cells.Item(23, 11).Value = 1.038751992887211
$ws.Cells.Item(23, 12).Value = 1.052602175019964
$ws.Cells.Item(23, 13).Value = 1.058394901599133
$ws.Cells.Item(23, 14).Value = 1.038062420514152

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.031064778025411
$ws.Cells.Item(24, 4).Value = 1.036028379200091
$ws.Cells.Item(24, 5).Value = 1.050557717406905
$ws.Cells.Item(24, 6).Value = 1.056395822828892
$ws.Cells.Item(24, 9).Value = 1.037636749775986
$ws.Cells.Item(24, 10).Value = 1.037287706187908
$ws.Cells.Item(24, 11).Value = 1.039398341599979
$ws.Cells.Item(24, 12).Value = 1.053877751270755
$ws.Cells.Item(24, 13).Value = 1.05969621952201
$ws.Cells.Item(24, 14).Value = 1.038760773709258

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.032482953986888
$ws.Cells.Item(25, 4).Value = 1.037098264004582
$ws.Cells.Item(25, 5).Value = 1.052359571089347
$ws.Cells.Item(25, 6).Value = 1.058225029032553
$ws.Cells.Item(25, 9).Value = 1.037978702160606
$ws.Cells.Item(25, 10).Value = 1.038095803624861
$ws.Cells.Item(25, 11).Value = 1.04014671879536
$ws.Cells.Item(25, 12).Value = 1.055360951862975
$ws.Cells.Item(25, 13).Value = 1.061208711333003
$ws.Cells.Item(25, 14).Value = 1.039570018737262
